# Apply the commit's changes to the "Монеты" (Coins) worksheet:
#  1. Re-label the price column header from "Цена" to "Цена, ориентировочно".
#  2. Update the USD->local currency exchange rate (O2) from 100 to 82 so the
#     "I" column (price = H * $O$2) recalculates.
#  3. Extend the AutoFilter / _FilterDatabase defined name from the header-only
#     range A2:E2 to the full data range A2:E53.
#  4. Move the active-cell selection on the data sheet to E60 (matches the
#     author's last on-screen selection after editing row 53 / col E).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Монеты")

# 1) Header relabel
$ws.Range("H1").Value = "Цена, ориентировочно"

# 2) Exchange-rate update -> ripples through column I via existing formulas
$ws.Range("O2").Value = 82

# 3) Re-apply AutoFilter over the full data range, and update the
#    hidden _FilterDatabase defined name to match.
$ws.AutoFilterMode = $false
$ws.Range("A2:E53").AutoFilter()

$filterName = $wb.Names.Item("_xlnm._FilterDatabase")
$filterName.RefersTo = "=Монеты!`$A`$2:`$E`$53"

# 4) Restore the author's final on-sheet selection
$ws.Activate()
$ws.Range("E60").Select()
